# Applies the "railway deployment setup and docs" update to the requirement
# status document. We work from the bottom of the document upward so that
# paragraph indices for not-yet-touched (earlier) paragraphs stay valid.

$d = $word.ActiveDocument

# --- Git state section (bottom of doc) -------------------------------------
# Para 26: "- Current document rename/update is local and not pushed yet."
$d.Paragraphs(26).Range.Text = "- Current Railway deployment updates are local and not pushed yet."

# Para 25: "- Last pushed commit: 1adb53a"
$d.Paragraphs(25).Range.Text = "- Last pushed commit: e839976"

# --- Insert new "Config" / "DB/Migrations" sections before "Git state" -----
# Para 23 is the blank line right before "Git state" (para 24); insert the
# new section's paragraphs right after that blank line.
$newSectionLines = @(
    "Config",
    "- Frontend Turnstile: REACT_APP_TURNSTILE_SITE_KEY",
    "- Backend Turnstile: Turnstile__SecretKey",
    "",
    "DB/Migrations",
    "- Ensure target DB has latest migrations: dotnet ef database update",
    ""
)
$insertAfter = $d.Paragraphs(23)
foreach ($line in $newSectionLines) {
    $insertAfter.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($insertAfter.Index + 1)
    $newPara.Range.Text = $line
    $insertAfter = $newPara
}

# --- Pending / Partial section ----------------------------------------------
# Para 20: drop "e.g., " and "editing" -> "editor"
$d.Paragraphs(20).Range.Text = "- Some advanced form fields are minimal for now (rich order item editor, full address management UI)."

# --- "Document structure update" -> "New in this update" section -----------
# Insert the 4th new bullet after para 16, then rewrite paras 13-16 in place.
$d.Paragraphs(16).Range.InsertParagraphAfter()
$d.Paragraphs(17).Range.Text = "  - ``ASPNETCORE_URLS=http://+:`${PORT:-5000}`` at runtime entrypoint."

$d.Paragraphs(16).Range.Text = "- Updated backend Dockerfile to bind Railway dynamic PORT using:"
$d.Paragraphs(15).Range.Text = "- Added frontend production Dockerfile: ``frontend/Dockerfile``."
$d.Paragraphs(14).Range.Text = "- Added Railway deployment runbook: ``DEPLOY_RAILWAY.md``."
$d.Paragraphs(13).Range.Text = "New in this update (Railway deploy readiness)"

# --- "Developed" section ----------------------------------------------------
# Para 11: drop "now"
$d.Paragraphs(11).Range.Text = "- Settings General save updates store via backend API."

# Para 10: delete entirely ("- Product/Customer/Order create+update+delete API wiring done end-to-end.")
$d.Paragraphs(10).Range.Delete()

# Para 6: append "(+ CRUD wiring in admin)"
$d.Paragraphs(6).Range.Text = "- Core domain endpoints: merchants, stores, products, customers, orders (+ CRUD wiring in admin)."

Write-Host "edit complete"
